$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay as text (they may look numeric,
# e.g. "320.38", or use "." as a thousands-style separator, e.g. "30.115.24").
# Force text format first so Excel does not reinterpret/round them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.115.24'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.916.91'
$ws.Range("E3").Value = '  +2.58%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.38'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5060'
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4083'
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08355'
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.44'
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.110'
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.02'
$ws.Range("E12").Value = '  +6.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.426'
$ws.Range("E13").Value = '  +2.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.909.56'
$ws.Range("E14").Value = '  +2.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.245'
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.58'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06517'
$ws.Range("E19").Value = '  +2.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.56'
$ws.Range("E20").Value = '  +4.11%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.947'
$ws.Range("E22").Value = '  +3.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.129.00'
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.35'
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.131.93'
$ws.Range("E26").Value = '  +2.57%  '
$ws.Range("E27").Value = '  +4.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.86'
$ws.Range("E28").Value = '  +1.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.278'
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.72'
$ws.Range("E30").Value = '  +1.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.146'
$ws.Range("E31").Value = '  +9.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1044'
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("E33").Value = '  +1.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.781'
$ws.Range("E34").Value = '  +1.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02460'
$ws.Range("E35").Value = '  +2.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.367'
$ws.Range("E36").Value = '  +3.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06438'
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2157'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6534'
$ws.Range("E39").Value = '  +4.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.199'
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.633'
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.42'
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.42'
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6088'
$ws.Range("E45").Value = '  +3.70%  '
$ws.Range("E46").Value = '  +10.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.622'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.211'
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.24'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '79.08'
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("E51").Value = '  -0.55%  '
